# Generate Report for Archive
# Update the localization status from "Ready for handoff" to "In Translation"
# across the Overview, zh-cn and de-de sheets, and shrink the now-narrower
# "Status" columns to match the report's refreshed auto-fit width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Replace the status text everywhere it appears.
$wsOverview.Cells.Replace("Ready for handoff", "In Translation")
$wsZhCn.Cells.Replace("Ready for handoff", "In Translation")
$wsDeDe.Cells.Replace("Ready for handoff", "In Translation")

# The status columns shrink now that "In Translation" is shorter than
# "Ready for handoff" - re-fit their widths.
$wsOverview.Columns.Item(5).ColumnWidth() = 12.5
$wsOverview.Columns.Item(6).ColumnWidth() = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth() = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth() = 12.5
